$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last daily-scrum block (rows 31-35) down to rows 37-41 to replicate formatting/styles
$ws.Range("A31:D35").Copy($ws.Range("A37"))

# Update the date for the new entry
$ws.Range("A37").Value2 = 44232

# Set values in an order that matches the desired shared-string allocation order
$ws.Range("C38").Value2 = "atualização e conclusao dos artefactos de análise"
$ws.Range("C39").Value2 = "atualização e conclusao dos artefactos de análise"
$ws.Range("C40").Value2 = "atualização e conclusao dos artefactos de análise"
$ws.Range("C41").Value2 = "atualização e conclusao dos artefactos de análise"

$ws.Range("B40").Value2 = "testes de utilização da aplicação, ajuste de pequenos bugs e funcionalidades`nupdate da lógica mvc"

$ws.Range("B38").Value2 = "testes de utilização da aplicação, ajuste de pequenos bugs e funcionalidades`njavadocs"
$ws.Range("B39").Value2 = "testes de utilização da aplicação, ajuste de pequenos bugs e funcionalidades`njavadocs"
$ws.Range("B41").Value2 = "testes de utilização da aplicação, ajuste de pequenos bugs e funcionalidades`njavadocs"

$ws.Range("D38").Value2 = "NADA A APONTAR"
$ws.Range("D39").Value2 = "NADA A APONTAR"
$ws.Range("D40").Value2 = "NADA A APONTAR"
$ws.Range("D41").Value2 = "NADA A APONTAR"

# Set the correct row heights for the new block
$ws.Rows(37).RowHeight = 162.45
$ws.Rows(38).RowHeight = 162.45
$ws.Rows(39).RowHeight = 162.45
$ws.Rows(40).RowHeight = 162.45
$ws.Rows(41).RowHeight = 162.45

# Update the viewport / selection to match final state
$ws.Range("B38").Select()
